$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in the title cell (A1)
$ws.Range("A1").Value = 'Datos actualizados a 31 de Marzo de 2020 a las 15:20'

# Refresh country data: a handful of rows picked up new case counts in this
# update, which also re-sorted a few countries relative to their neighbours
# (the table is ordered by total cases, descending).

# Row 16: Austria
$ws.Cells.Item(16, 1).Value = 'Austria'
$ws.Cells.Item(16, 2).Value = 9974
$ws.Cells.Item(16, 3).Value = 356
$ws.Cells.Item(16, 4).Value = 1095
$ws.Cells.Item(16, 5).Value = 8751
$ws.Cells.Item(16, 6).Value = 198
$ws.Cells.Item(16, 7).Value = 20
$ws.Cells.Item(16, 8).Value = 128

# Row 39: Arabia Saudita
$ws.Cells.Item(39, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(39, 2).Value = 1563
$ws.Cells.Item(39, 3).Value = 110
$ws.Cells.Item(39, 4).Value = 165
$ws.Cells.Item(39, 5).Value = 1388
$ws.Cells.Item(39, 6).Value = 12
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = 10

# Row 48: Argentina
$ws.Cells.Item(48, 1).Value = 'Argentina'
$ws.Cells.Item(48, 2).Value = 966
$ws.Cells.Item(48, 3).Value = 146
$ws.Cells.Item(48, 4).Value = 228
$ws.Cells.Item(48, 5).Value = 713
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 2
$ws.Cells.Item(48, 8).Value = 25

# Row 52: Serbia
$ws.Cells.Item(52, 1).Value = 'Serbia'
$ws.Cells.Item(52, 2).Value = 900
$ws.Cells.Item(52, 3).Value = 115
$ws.Cells.Item(52, 4).Value = 42
$ws.Cells.Item(52, 5).Value = 835
$ws.Cells.Item(52, 6).Value = 62
$ws.Cells.Item(52, 7).Value = 7
$ws.Cells.Item(52, 8).Value = 23

# Row 53: Croacia
$ws.Cells.Item(53, 1).Value = 'Croacia'
$ws.Cells.Item(53, 2).Value = 867
$ws.Cells.Item(53, 3).Value = 77
$ws.Cells.Item(53, 4).Value = 67
$ws.Cells.Item(53, 5).Value = 794
$ws.Cells.Item(53, 6).Value = 32
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 6

# Row 54: Eslovenia
$ws.Cells.Item(54, 1).Value = 'Eslovenia'
$ws.Cells.Item(54, 2).Value = 802
$ws.Cells.Item(54, 3).Value = 46
$ws.Cells.Item(54, 4).Value = 10
$ws.Cells.Item(54, 5).Value = 777
$ws.Cells.Item(54, 6).Value = 24
$ws.Cells.Item(54, 7).Value = 4
$ws.Cells.Item(54, 8).Value = 15

# Row 55: Colombia
$ws.Cells.Item(55, 1).Value = 'Colombia'
$ws.Cells.Item(55, 2).Value = 798
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 15
$ws.Cells.Item(55, 5).Value = 769
$ws.Cells.Item(55, 6).Value = 29
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 14

# Row 80: Republica de Macedonia
$ws.Cells.Item(80, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(80, 2).Value = 329
$ws.Cells.Item(80, 3).Value = 44
$ws.Cells.Item(80, 4).Value = 12
$ws.Cells.Item(80, 5).Value = 310
$ws.Cells.Item(80, 6).Value = 1
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 7

# Row 81: Taiwan
$ws.Cells.Item(81, 1).Value = 'Taiwan'
$ws.Cells.Item(81, 2).Value = 322
$ws.Cells.Item(81, 3).Value = 16
$ws.Cells.Item(81, 4).Value = 39
$ws.Cells.Item(81, 5).Value = 278
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 5

# Row 82: Uruguay
$ws.Cells.Item(82, 1).Value = 'Uruguay'
$ws.Cells.Item(82, 2).Value = 320
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 25
$ws.Cells.Item(82, 5).Value = 294
$ws.Cells.Item(82, 6).Value = 9
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 1

# Row 83: Moldavia
$ws.Cells.Item(83, 1).Value = 'Moldavia'
$ws.Cells.Item(83, 2).Value = 298
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 18
$ws.Cells.Item(83, 5).Value = 278
$ws.Cells.Item(83, 6).Value = 44
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 2

# Row 84: Azerbaiyan
$ws.Cells.Item(84, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(84, 2).Value = 298
$ws.Cells.Item(84, 3).Value = 25
$ws.Cells.Item(84, 4).Value = 26
$ws.Cells.Item(84, 5).Value = 267
$ws.Cells.Item(84, 6).Value = 11
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 5

# Row 85: Kuwait
$ws.Cells.Item(85, 1).Value = 'Kuwait'
$ws.Cells.Item(85, 2).Value = 289
$ws.Cells.Item(85, 3).Value = 23
$ws.Cells.Item(85, 4).Value = 73
$ws.Cells.Item(85, 5).Value = 216
$ws.Cells.Item(85, 6).Value = 13
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0

# Row 158: Birmania
$ws.Cells.Item(158, 1).Value = 'Birmania'
$ws.Cells.Item(158, 2).Value = 14
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 13
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 1
$ws.Cells.Item(158, 8).Value = 1

# Row 159: Bahamas
$ws.Cells.Item(159, 1).Value = 'Bahamas'
$ws.Cells.Item(159, 2).Value = 14
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 1
$ws.Cells.Item(159, 5).Value = 13
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

# Row 169: Granada
$ws.Cells.Item(169, 1).Value = 'Granada'
$ws.Cells.Item(169, 2).Value = 9
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

# Row 170: Laos
$ws.Cells.Item(170, 1).Value = 'Laos'
$ws.Cells.Item(170, 2).Value = 9
$ws.Cells.Item(170, 3).Value = 1
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 9
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 0

# Row 171: Suazilandia
$ws.Cells.Item(171, 1).Value = 'Suazilandia'
$ws.Cells.Item(171, 2).Value = 9
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 9
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0

# Row 173: San Cristobal y Nieves
$ws.Cells.Item(173, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(173, 2).Value = 8
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 8
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

# Row 174: Surinam
$ws.Cells.Item(174, 1).Value = 'Surinam'
$ws.Cells.Item(174, 2).Value = 8
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 8
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

# Row 176: Mozambique
$ws.Cells.Item(176, 1).Value = 'Mozambique'
$ws.Cells.Item(176, 2).Value = 8
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 8
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

# Row 177: Guinea-Bisau
$ws.Cells.Item(177, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(177, 2).Value = 8
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 8
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

# Row 180: Antigua y Barbuda
$ws.Cells.Item(180, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(180, 2).Value = 7
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

# Row 181: Republica del Chad
$ws.Cells.Item(181, 1).Value = 'Republica del Chad'
$ws.Cells.Item(181, 2).Value = 7
$ws.Cells.Item(181, 3).Value = 2
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 7
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

# Row 186: Benin
$ws.Cells.Item(186, 1).Value = 'Benin'
$ws.Cells.Item(186, 2).Value = 6
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 1
$ws.Cells.Item(186, 5).Value = 5
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 187: Cabo Verde
$ws.Cells.Item(187, 1).Value = 'Cabo Verde'
$ws.Cells.Item(187, 2).Value = 6
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 5
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 1

# Row 188: San Bartolome
$ws.Cells.Item(188, 1).Value = 'San Bartolome'
$ws.Cells.Item(188, 2).Value = 6
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 1
$ws.Cells.Item(188, 5).Value = 5
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

# Row 190: Fiyi
$ws.Cells.Item(190, 1).Value = 'Fiyi'
$ws.Cells.Item(190, 2).Value = 5
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# Row 191: Islas Turcas y Caicos
$ws.Cells.Item(191, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(191, 2).Value = 5
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 0

# Row 192: Montserrat
$ws.Cells.Item(192, 1).Value = 'Montserrat'
$ws.Cells.Item(192, 2).Value = 5
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

# Row 197: Botsuana
$ws.Cells.Item(197, 1).Value = 'Botsuana'
$ws.Cells.Item(197, 2).Value = 3
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 3
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0

# Row 198: Islas Virgenes Britanicas
$ws.Cells.Item(198, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(198, 2).Value = 3
$ws.Cells.Item(198, 3).Value = 1
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 5).Value = 3
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0

# Row 199: Belice
$ws.Cells.Item(199, 1).Value = 'Belice'
$ws.Cells.Item(199, 2).Value = 3
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 3
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# Row 200: Republica de Africa Central
$ws.Cells.Item(200, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(200, 2).Value = 3
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 3
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0

# Row 201: Liberia
$ws.Cells.Item(201, 1).Value = 'Liberia'
$ws.Cells.Item(201, 2).Value = 3
$ws.Cells.Item(201, 3).Value = 0
$ws.Cells.Item(201, 4).Value = 0
$ws.Cells.Item(201, 5).Value = 3
$ws.Cells.Item(201, 6).Value = 0
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 0
